# Daily attendance processing - swap the order of the first two
# comma-separated "Recorded By" names in column G for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
